$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff (cryptos list refresh, Wed Sep 20 19:08:46 UTC 2023)
# NumberFormat "@" forces values that look numeric (e.g. "4.41", "216.29")
# to be stored as text, matching the original inlineStr cells; resetting the
# Style back to "Normal" afterwards avoids leaving a stray quote-prefix style
# on the cell so only the cell VALUES change, exactly like the source diff.
$updates = @(
    @{ Cell = "D2"; Value = "27.204.40" }
    @{ Cell = "E2"; Value = "  +0.05%  " }
    @{ Cell = "D3"; Value = "1.633.81" }
    @{ Cell = "E3"; Value = "  -0.85%  " }
    @{ Cell = "E4"; Value = "  -0.01%  " }
    @{ Cell = "D5"; Value = "216.29" }
    @{ Cell = "E5"; Value = "  -0.44%  " }
    @{ Cell = "E6"; Value = "  +1.41%  " }
    @{ Cell = "E8"; Value = "  +0.00%  " }
    @{ Cell = "E9"; Value = "  -0.16%  " }
    @{ Cell = "D10"; Value = "20.29" }
    @{ Cell = "E10"; Value = "  +1.74%  " }
    @{ Cell = "D11"; Value = "0.0849" }
    @{ Cell = "E11"; Value = "  +0.00%  " }
    @{ Cell = "D12"; Value = "1.636.30" }
    @{ Cell = "E12"; Value = "  -0.68%  " }
    @{ Cell = "E13"; Value = "  +0.13%  " }
    @{ Cell = "D14"; Value = "0.545" }
    @{ Cell = "E14"; Value = "  +0.76%  " }
    @{ Cell = "D15"; Value = "65.18" }
    @{ Cell = "E15"; Value = "  -3.59%  " }
    @{ Cell = "D16"; Value = "27.196.72" }
    @{ Cell = "E16"; Value = "  +0.06%  " }
    @{ Cell = "D17"; Value = "0.0₃0742" }
    @{ Cell = "E17"; Value = "  +0.31%  " }
    @{ Cell = "D18"; Value = "217.94" }
    @{ Cell = "E18"; Value = "  -0.55%  " }
    @{ Cell = "E19"; Value = "  -0.03%  " }
    @{ Cell = "E20"; Value = "  +1.62%  " }
    @{ Cell = "D21"; Value = "4.41" }
    @{ Cell = "E21"; Value = "  -0.21%  " }
    @{ Cell = "E22"; Value = "  -6.50%  " }
    @{ Cell = "D23"; Value = "9.08" }
    @{ Cell = "E23"; Value = "  -1.56%  " }
    @{ Cell = "D24"; Value = "148.15" }
    @{ Cell = "E24"; Value = "  +0.33%  " }
    @{ Cell = "E25"; Value = "  +0.03%  " }
    @{ Cell = "E26"; Value = "  -3.00%  " }
    @{ Cell = "E27"; Value = "  +0.50%  " }
    @{ Cell = "D28"; Value = "15.69" }
    @{ Cell = "E28"; Value = "  -0.48%  " }
    @{ Cell = "E29"; Value = "  -0.09%  " }
    @{ Cell = "E30"; Value = "  -0.35%  " }
    @{ Cell = "E31"; Value = "  -0.48%  " }
    @{ Cell = "E32"; Value = "  -1.01%  " }
    @{ Cell = "D33"; Value = "1.343.89" }
    @{ Cell = "E33"; Value = "  +6.13%  " }
    @{ Cell = "E34"; Value = "  +0.13%  " }
    @{ Cell = "E35"; Value = "  -0.10%  " }
    @{ Cell = "E36"; Value = "  -0.47%  " }
    @{ Cell = "D37"; Value = "0.548" }
    @{ Cell = "E38"; Value = "  +0.33%  " }
    @{ Cell = "E39"; Value = "  -0.06%  " }
    @{ Cell = "E40"; Value = "  +1.28%  " }
    @{ Cell = "D41"; Value = "0.803" }
    @{ Cell = "E41"; Value = "  -0.64%  " }
    @{ Cell = "D42"; Value = "64.78" }
    @{ Cell = "E42"; Value = "  +4.50%  " }
    @{ Cell = "D43"; Value = "5.28" }
    @{ Cell = "E43"; Value = "  -3.04%  " }
    @{ Cell = "D44"; Value = "1.773.36" }
    @{ Cell = "E44"; Value = "  -0.92%  " }
    @{ Cell = "D45"; Value = "90.85" }
    @{ Cell = "E45"; Value = "  -0.94%  " }
    @{ Cell = "D46"; Value = "1.61" }
    @{ Cell = "E46"; Value = "  +0.18%  " }
    @{ Cell = "D47"; Value = "0.812" }
    @{ Cell = "E47"; Value = "  +21.91%  " }
    @{ Cell = "B48"; Value = "Cronos" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" }
    @{ Cell = "D48"; Value = "0.0514" }
    @{ Cell = "E48"; Value = "  +0.06%  " }
    @{ Cell = "B49"; Value = "Algorand" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" }
    @{ Cell = "D49"; Value = "0.0990" }
    @{ Cell = "E49"; Value = "  +1.70%  " }
    @{ Cell = "B50"; Value = "EnergySwap" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D50"; Value = "7.60" }
    @{ Cell = "E50"; Value = "  -0.46%  " }
    @{ Cell = "B51"; Value = "USDD" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd" }
    @{ Cell = "D51"; Value = "1.00" }
    @{ Cell = "E51"; Value = "  -0.34%  " }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}